$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update 想去人数 (attendee counts) in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5392
$wsExhibit.Range("F4").Value = 11459
$wsExhibit.Range("F5").Value = 282
$wsExhibit.Range("F7").Value = 168
$wsExhibit.Range("F8").Value = 258
$wsExhibit.Range("F9").Value = 987
$wsExhibit.Range("F10").Value = 96

# Sheet "全部类型" (All types) - same updates, different row positions
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5392
$wsAll.Range("F7").Value = 11459
$wsAll.Range("F8").Value = 282
$wsAll.Range("F10").Value = 168
$wsAll.Range("F13").Value = 258
$wsAll.Range("F14").Value = 987
$wsAll.Range("F16").Value = 96
